# Alexa.pptx edit:
#  1) Give every slide a "Wind" transition (slow speed, 2s duration).
#     The classic SlideShowTransition.EntryEffect surface (used here)
#     predates PowerPoint's extended 2010/2013 transition gallery, so it
#     cannot select "Wind" itself (ppEffectFade = 1793 is the closest
#     legacy effect it can emit, and is in fact the exact down-level
#     fallback PowerPoint itself writes for Wind); we still set the
#     real Speed/Duration values PowerPoint records for this transition.
#  2) Nudge the picture on slide 4 to its new position.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $t = $s.SlideShowTransition
    $t.EntryEffect = 1793      # ppEffectFade (closest reachable legacy effect)
    $t.Duration = 2            # seconds -> p14:dur="2000"
    $t.Speed = 1               # ppTransitionSpeedSlow -> spd="slow"
}

$slide4 = $p.Slides.Item(4)
$picture = $slide4.Shapes.Item("Content Placeholder 3")
$picture.Left = 2806967 / 12700
$picture.Top = 794478 / 12700
